{"js": "// Each entry is [oldEquationText, newEquationText]. All 25 equation cells\n// in the single table on the page are updated from the old values to the\n// new ones, matching the commit's regenerated numbers.\nconst pairs = [\n  [\"65\u00d714=910\", \"45\u00d793=4185\"],\n  [\"75\u00d744=3300\", \"39\u00d779=3081\"],\n  [\"78\u00d777=6006\", \"34\u00d780=2720\"],\n  [\"73\u00d748=3504\", \"27\u00d785=2295\"],\n  [\"63\u00d790=5670\", \"96\u00d766=6336\"],\n  [\"93\u00d754=5022\", \"45\u00d723=1035\"],\n  [\"89\u00d733=2937\", \"79\u00d789=7031\"],\n  [\"60\u00d787=5220\", \"60\u00d798=5880\"],\n  [\"42\u00d753=2226\", \"79\u00d790=7110\"],\n  [\"15\u00d775=1125\", \"61\u00d745=2745\"],\n  [\"43\u00d722=946\", \"55\u00d735=1925\"],\n  [\"61\u00d773=4453\", \"51\u00d768=3468\"],\n  [\"64\u00d713=832\", \"88\u00d724=2112\"],\n  [\"54\u00d795=5130\", \"79\u00d716=1264\"],\n  [\"51\u00d773=3723\", \"33\u00d718=594\"],\n  [\"56\u00d722=1232\", \"66\u00d763=4158\"],\n  [\"90\u00d721=1890\", \"12\u00d770=840\"],\n  [\"84\u00d721=1764\", \"29\u00d787=2523\"],\n  [\"13\u00d728=364\", \"77\u00d756=4312\"],\n  [\"74\u00d755=4070\", \"25\u00d715=375\"],\n  [\"29\u00d771=2059\", \"17\u00d760=1020\"],\n  [\"98\u00d757=5586\", \"41\u00d784=3444\"],\n  [\"21\u00d731=651\", \"58\u00d771=4118\"],\n  [\"76\u00d768=5168\", \"72\u00d732=2304\"],\n  [\"85\u00d788=7480\", \"78\u00d784=6552\"],\n];\n\nfor (const [oldText, newText] of pairs) {\n  // Find the unique run containing the old equation text and replace it\n  // in place so the run's formatting (font, size, ...) is preserved.\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length !== 1) {\n    throw new Error(`Expected exactly 1 match for \"${oldText}\", found ${results.items.length}`);\n  }\n\n  results.items[0].insertText(newText, \"Replace\");\n}\n\nawait context.sync();\n", "ps1": "# Each pair is (oldEquationText, newEquationText). All 25 equation cells\n# in the document's single table are updated from the old values to the\n# new ones, matching the commit's regenerated numbers.\n$d = $word.ActiveDocument\n\n$pairs = @(\n  @(\"65\u00d714=910\", \"45\u00d793=4185\"),\n  @(\"75\u00d744=3300\", \"39\u00d779=3081\"),\n  @(\"78\u00d777=6006\", \"34\u00d780=2720\"),\n  @(\"73\u00d748=3504\", \"27\u00d785=2295\"),\n  @(\"63\u00d790=5670\", \"96\u00d766=6336\"),\n  @(\"93\u00d754=5022\", \"45\u00d723=1035\"),\n  @(\"89\u00d733=2937\", \"79\u00d789=7031\"),\n  @(\"60\u00d787=5220\", \"60\u00d798=5880\"),\n  @(\"42\u00d753=2226\", \"79\u00d790=7110\"),\n  @(\"15\u00d775=1125\", \"61\u00d745=2745\"),\n  @(\"43\u00d722=946\", \"55\u00d735=1925\"),\n  @(\"61\u00d773=4453\", \"51\u00d768=3468\"),\n  @(\"64\u00d713=832\", \"88\u00d724=2112\"),\n  @(\"54\u00d795=5130\", \"79\u00d716=1264\"),\n  @(\"51\u00d773=3723\", \"33\u00d718=594\"),\n  @(\"56\u00d722=1232\", \"66\u00d763=4158\"),\n  @(\"90\u00d721=1890\", \"12\u00d770=840\"),\n  @(\"84\u00d721=1764\", \"29\u00d787=2523\"),\n  @(\"13\u00d728=364\", \"77\u00d756=4312\"),\n  @(\"74\u00d755=4070\", \"25\u00d715=375\"),\n  @(\"29\u00d771=2059\", \"17\u00d760=1020\"),\n  @(\"98\u00d757=5586\", \"41\u00d784=3444\"),\n  @(\"21\u00d731=651\", \"58\u00d771=4118\"),\n  @(\"76\u00d768=5168\", \"72\u00d732=2304\"),\n  @(\"85\u00d788=7480\", \"78\u00d784=6552\")\n)\n\nforeach ($pair in $pairs) {\n  $oldText = $pair[0]\n  $newText = $pair[1]\n\n  # Find the unique occurrence of the old equation text and replace it in\n  # place so the run's formatting (font, size, ...) is preserved.\n  $find = $d.Content.Find\n  $find.ClearFormatting()\n  $find.Replacement.ClearFormatting()\n  $find.Text = $oldText\n  $find.Replacement.Text = $newText\n  $found = $find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n\n  if (-not $found) {\n    throw \"Could not find expected text: $oldText\"\n  }\n}\n\n"}
